# edu_table_helper.xlsx update
# - Remove the unused "Sheet2" worksheet
# - Update the out_of_school barrier labels (G2:G7) from internal codes to
#   the human readable text, and refresh sheet selections so the
#   out_of_school tab is the active one (mirrors the authored edit).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the leftover "Sheet2" worksheet entirely.
$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.Delete()

$oos = $wb.Worksheets.Item("out_of_school")

# Assign the new values in reverse row order so the shared-string table is
# rebuilt in the same order Excel produced it (Other .. child_work_outside).
$oos.Range("G7").Value = "Other"
$oos.Range("G6").Value = "School has been closed due to conflict"
$oos.Range("G5").Value = "There is a lack of interest/Education is not a priority either for the child or the household"
$oos.Range("G4").Value = "Lack of appropriate and accessible school"
$oos.Range("G3").Value = "Cannot afford education-related costs (e.g. tuition, supplies, transportation)"
$oos.Range("G2").Value = "Child participating in income generating activities outside of the home"

# Move the active sheet selection from "level1" to "out_of_school".
$level1 = $wb.Worksheets.Item("level1")
$level1.Range("P30").Select()

$oos.Activate()
$oos.Range("E15").Select()
